$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.365740299224854
$ws.Range("B1").Value = 1.674098253250122
$ws.Range("C1").Value = 2.325392246246338
$ws.Range("D1").Value = 4.989462852478027
$ws.Range("E1").Value = 1.871124267578125
